$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.819.60"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "3.527.15"
$ws.Range("E3").Value = "  -3.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "196.28"
$ws.Range("E5").Value = "  -2.63%  "

# Row 6
$ws.Range("D6").Value = "554.47"
$ws.Range("E6").Value = "  -2.85%  "

# Row 7
$ws.Range("D7").Value = "0.655"
$ws.Range("E7").Value = "  +5.75%  "

# Row 8
$ws.Range("D8").Value = "3.522.64"
$ws.Range("E8").Value = "  -2.70%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").Value = "0.662"
$ws.Range("E10").Value = "  -2.69%  "

# Row 11
$ws.Range("D11").Value = "60.31"
$ws.Range("E11").Value = "  +3.63%  "

# Row 12
$ws.Range("D12").Value = "0.144"
$ws.Range("E12").Value = "  -6.49%  "

# Row 13
$ws.Range("D13").Value = "0.0000268"
$ws.Range("E13").Value = "  -8.58%  "

# Row 14
$ws.Range("D14").Value = "9.91"
$ws.Range("E14").Value = "  -1.54%  "

# Row 15
$ws.Range("D15").Value = "4.104.34"
$ws.Range("E15").Value = "  -2.58%  "

# Row 16
$ws.Range("D16").Value = "3.536.57"
$ws.Range("E16").Value = "  -2.78%  "

# Row 17
$ws.Range("E17").Value = "  -1.68%  "

# Row 18
$ws.Range("D18").Value = "67.671.23"
$ws.Range("E18").Value = "  -0.82%  "

# Row 19
$ws.Range("D19").Value = "18.33"
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("D20").Value = "11.87"
$ws.Range("E20").Value = "  -4.49%  "

# Row 21
$ws.Range("E21").Value = "  -5.33%  "

# Row 22
$ws.Range("D22").Value = "401.76"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").Value = "87.27"
$ws.Range("E23").Value = "  +1.80%  "

# Row 24
$ws.Range("D24").Value = "3.98"
$ws.Range("E24").Value = "  -5.60%  "

# Row 25
$ws.Range("D25").Value = "11.66"
$ws.Range("E25").Value = "  -11.21%  "

# Row 26
$ws.Range("D26").Value = "12.36"
$ws.Range("E26").Value = "  -1.70%  "

# Row 27
$ws.Range("D27").Value = "2.83"
$ws.Range("E27").Value = "  -4.56%  "

# Row 28
$ws.Range("D28").Value = "3.83"
$ws.Range("E28").Value = "  -0.50%  "

# Row 29
$ws.Range("D29").Value = "8.88"
$ws.Range("E29").Value = "  -3.20%  "

# Row 30
$ws.Range("D30").Value = "716.67"
$ws.Range("E30").Value = "  +4.06%  "

# Row 31
$ws.Range("D31").Value = "31.40"
$ws.Range("E31").Value = "  -1.65%  "

# Row 32
$ws.Range("D32").Value = "7.01"
$ws.Range("E32").Value = "  -13.89%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  -4.19%  "

# Row 34
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "64.30"
$ws.Range("E34").Value = "  -0.09%  "

# Row 35
$ws.Range("E35").Value = "  -4.15%  "

# Row 36
$ws.Range("D36").Value = "38.41"
$ws.Range("E36").Value = "  -10.25%  "

# Row 37
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").Value = "0.390"
$ws.Range("E38").Value = "  -8.50%  "

# Row 39
$ws.Range("E39").Value = "  -5.06%  "

# Row 40
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.056.68"
$ws.Range("E41").Value = "  -6.46%  "

# Row 42
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").Value = "2.99"
$ws.Range("E42").Value = "  -4.87%  "

# Row 43
$ws.Range("D43").Value = "0.0₃0680"
$ws.Range("E43").Value = "  -13.12%  "

# Row 44
$ws.Range("D44").Value = "2.77"
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$ws.Range("D45").Value = "0.136"
$ws.Range("E45").Value = "  +2.63%  "

# Row 46
$ws.Range("D46").Value = "2.48"
$ws.Range("E46").Value = "  -12.05%  "

# Row 47
$ws.Range("D47").Value = "0.0409"
$ws.Range("E47").Value = "  -2.58%  "

# Row 48
$ws.Range("D48").Value = "139.57"
$ws.Range("E48").Value = "  -0.69%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.01"
$ws.Range("E49").Value = "  -2.98%  "

# Row 50
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").Value = "  -16.34%  "

# Row 51
$ws.Range("D51").Value = "8.28"
$ws.Range("E51").Value = "  -7.23%  "
